$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Sheet1 -> Worksheet1)
$ws.Name = "Worksheet1"

# Resize the data columns (A, B, C) to their new widths
$ws.Columns.Item(1).ColumnWidth = 30.89115646258507
$ws.Columns.Item(2).ColumnWidth = 18.47278911564627
$ws.Columns.Item(3).ColumnWidth = 22.38605442176867

# Tighten the two data rows back to the sheet's normal row height
$ws.Rows.Item(1).RowHeight = 12.8
$ws.Rows.Item(2).RowHeight = 12.8

# Move the active selection to B6
$ws.Range("B6").Select() | Out-Null

# Update page margins (now using Normal-style margins) and clear the
# custom header/footer text
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36.85
$ps.FooterMargin = 36.85
$ps.CenterHeader = ""
$ps.CenterFooter = ""
